# Regenerate merged AHB files
#
# The sheet "AHB-Diff" compares two AHB (Anwendungshandbuch) format versions.
# The header row previously labeled the two compared sides as "_old"/"_new";
# they are renamed here to the concrete format versions "_FV2410"/"_FV2504".
# In addition, the data range is turned into a real Excel Table, and the
# header row is frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "<name>_old" -> "<name>_FV2410" (columns A-J) ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

# Column K ("diff") is unchanged.

# --- Rename header row: "<name>_new" -> "<name>_FV2504" (columns L-U) ---
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- Turn the whole data range into an Excel Table (Table1) ---
$dataRange = $ws.Range("A1:U76")
$listObject = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$listObject.Name = "Table1"

# --- Freeze the header row (split/freeze above row 2) ---
$ws.Range("A2").Select() | Out-Null
[void]($ws.Application.ActiveWindow.FreezePanes = $true)
$ws.Range("A1").Select() | Out-Null
